# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# sheet with refreshed values, per the GitHub Actions data refresh commit.
#
# D-column values are written with a leading apostrophe (forces text,
# preventing Excel from re-interpreting number-like strings such as
# "1.00" or "0.160" as numeric values and dropping trailing zeros /
# separators), and the cell Style is then reset back to "Normal" so no
# stray quote-prefix style is left on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.Value = "'58.532.91"
$c.Style = 'Normal'
$ws.Range('E2').Value = '  -3.84%  '
$c = $ws.Range('D3')
$c.Value = "'2.539.76"
$c.Style = 'Normal'
$ws.Range('E3').Value = '  -3.60%  '
$ws.Range('E4').Value = '  -0.08%  '
$c = $ws.Range('D5')
$c.Value = "'507.79"
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -3.95%  '
$c = $ws.Range('D6')
$c.Value = "'144.13"
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -7.06%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -4.25%  '
$c = $ws.Range('D9')
$c.Value = "'2.543.16"
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -3.88%  '
$c = $ws.Range('D10')
$c.Value = "'6.09"
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -8.68%  '
$ws.Range('E11').Value = '  -6.24%  '
$ws.Range('E12').Value = '  -5.42%  '
$ws.Range('E13').Value = '  -0.60%  '
$c = $ws.Range('D14')
$c.Value = "'2.983.38"
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -3.61%  '
$c = $ws.Range('D15')
$c.Value = "'58.490.35"
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -3.94%  '
$c = $ws.Range('D16')
$c.Value = "'20.68"
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -5.91%  '
$ws.Range('E17').Value = '  -6.03%  '
$c = $ws.Range('D18')
$c.Value = "'2.538.93"
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -3.75%  '
$ws.Range('E19').Value = '  -4.81%  '
$c = $ws.Range('D20')
$c.Value = "'335.08"
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -5.16%  '
$c = $ws.Range('D21')
$c.Value = "'10.09"
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -5.02%  '
$ws.Range('E22').Value = '  -0.26%  '
$c = $ws.Range('D23')
$c.Value = "'5.96"
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -4.46%  '
$c = $ws.Range('D24')
$c.Value = "'60.61"
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -1.71%  '
$ws.Range('E25').Value = '  -4.71%  '
$c = $ws.Range('D26')
$c.Value = "'0.999"
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -0.11%  '
$c = $ws.Range('D27')
$c.Value = "'0.160"
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -4.98%  '
$c = $ws.Range('D28')
$c.Value = "'2.650.77"
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -3.53%  '
$ws.Range('E29').Value = '  -9.18%  '
$c = $ws.Range('D30')
$c.Value = "'6.97"
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -5.90%  '
$c = $ws.Range('D31')
$c.Value = "'1.00"
$c.Style = 'Normal'
$c = $ws.Range('D32')
$c.Value = "'149.68"
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -0.49%  '
$c = $ws.Range('D33')
$c.Value = "'5.85"
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -4.88%  '
$c = $ws.Range('D34')
$c.Value = "'18.54"
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -4.87%  '
$c = $ws.Range('D35')
$c.Value = "'1.54"
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -5.17%  '
$c = $ws.Range('D36')
$c.Value = "'0.920"
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +3.90%  '
$c = $ws.Range('D37')
$c.Value = "'3.91"
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -6.15%  '
$ws.Range('E38').Value = '  -7.48%  '
$c = $ws.Range('D39')
$c.Value = "'36.01"
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -1.68%  '
$ws.Range('E40').Value = '  -11.29%  '
$ws.Range('E41').Value = '  -6.92%  '
$c = $ws.Range('D42')
$c.Value = "'283.71"
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -7.27%  '
$ws.Range('E43').Value = '  -7.83%  '
$ws.Range('E44').Value = '  -2.55%  '
$c = $ws.Range('D45')
$c.Value = "'0.997"
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -0.12%  '
$c = $ws.Range('D46')
$c.Value = "'0.599"
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -6.59%  '
$ws.Range('E47').Value = '  -5.09%  '
$c = $ws.Range('D48')
$c.Value = "'18.66"
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -5.48%  '
$c = $ws.Range('D49')
$c.Value = "'10.28"
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -0.58%  '
$ws.Range('E50').Value = '  -5.08%  '
$ws.Range('E51').Value = '  -8.24%  '
